# Feature Toggle implementation + conclusao busca_ai: refresh transaction log with new entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=2; A=45235; B="Point Da Irae"; C="-60,00"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=3; A=45235; B="Boteco Parô"; C="-343,39"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=4; A=45235; B="Cuor Di Crema"; C="-19,50"; D="CINTHIA"; E="1/1"; F="" },
  @{ Row=5; A=45235; B="Emporium Sao Paulo"; C="-15,00"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=6; A=45234; B="Uber"; C="-29,97"; D="PHILIPPE"; E="1/1"; F="TRANSPORTE" },
  @{ Row=7; A=45234; B="Top"; C="-8,80"; D="PHILIPPE"; E="1/1"; F="TRANSPORTE" },
  @{ Row=8; A=45234; B="Motorsport"; C="-300,00"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=9; A=45234; B="Caldo De Cana Da"; C="-8,00"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=10; A=45234; B="Alem Do Hamburguer"; C="-20,00"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=11; A=45234; B="Top"; C="-17,60"; D="PHILIPPE"; E="1/1"; F="TRANSPORTE" },
  @{ Row=12; A=45233; B="Saru Sushi"; C="-650,87"; D="CINTHIA"; E="1/1"; F="" },
  @{ Row=13; A=45233; B="Bar Jobim"; C="-48,18"; D="PHILIPPE"; E="1/1"; F="LAZER" },
  @{ Row=14; A=45233; B="Bar Jobim"; C="-100,00"; D="PHILIPPE"; E="1/1"; F="LAZER" },
  @{ Row=15; A=45233; B="Centauro"; C="-179,99"; D="CINTHIA"; E="1/2"; F="" },
  @{ Row=16; A=45233; B="Panvel"; C="-183,80"; D="CINTHIA"; E="1/2"; F="" },
  @{ Row=17; A=45233; B="Panvel"; C="-0,00"; D="CINTHIA"; E="1/1"; F="" },
  @{ Row=18; A=45233; B="Oba Hortifruti"; C="-82,04"; D="PHILIPPE"; E="1/1"; F="MERCADO" },
  @{ Row=19; A=45233; B="Swift"; C="-55,96"; D="PHILIPPE"; E="1/1"; F="MERCADO" },
  @{ Row=20; A=45233; B="Senhora Adega"; C="-45,99"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=21; A=45233; B="Echope"; C="-179,64"; D="PHILIPPE"; E="1/1"; F="LAZER" },
  @{ Row=22; A=45232; B="Bar Jobim"; C="-90,00"; D="CINTHIA"; E="1/1"; F="LAZER" },
  @{ Row=23; A=45232; B="Bar Jobim"; C="-355,57"; D="PHILIPPE"; E="1/1"; F="LAZER" },
  @{ Row=24; A=45232; B="Koa Moema"; C="-49,80"; D="CINTHIA"; E="1/1"; F="ALIMENTAÇÃO" },
  @{ Row=25; A=45232; B="Veneza Enxovais Textil"; C="-305,52"; D="CINTHIA"; E="1/2"; F="" },
  @{ Row=26; A=45231; B="Sancto Churrasco"; C="-201,08"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=27; A=45231; B="Aga"; C="-8,00"; D="CINTHIA"; E="1/1"; F="" },
  @{ Row=28; A=45231; B="Abastece Ai"; C="-225,40"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=29; A=45231; B="Enel"; C="-225,71"; D="PHILIPPE"; E="1/1"; F="" },
  @{ Row=30; A=45231; B="Mercado Livre"; C="-67,36"; D="CINTHIA"; E="1/1"; F="" },
  @{ Row=31; A=45231; B="Mercado Livre"; C="-183,04"; D="CINTHIA"; E="1/1"; F="" }
)

# Clear any stale CATEGORIA values in the range we are about to rewrite,
# since row positions are being reshuffled.
$ws.Range("F2:F31").ClearContents()

foreach ($row in $rows) {
  $r = $row.Row
  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  if ($row.F -ne "") {
    $ws.Cells.Item($r, 6).Value = $row.F
  }
}
